# Insert two new rows (403 and 404) above the current row 403 in the
# "Femacal de La Calera - Zapallo" consolidated sheet, pushing the
# existing rows 403-428 down to 405-430, then populate the two new
# rows with the new weekly price-survey records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows before the current row 403 (shifts 403:428 -> 405:430)
$ws.Rows("403:404").Insert()

# --- New row 403: Camote, 1a (cosecha) ---
$ws.Range("A403").Value = 3
$ws.Range("B403").Value = "Femacal de La Calera"
$ws.Range("C403").Value = "Coquimbo"
$ws.Range("D403").Value = 44585
$ws.Range("E403").Value = 5
$ws.Range("F403").Value = 100112045
$ws.Range("G403").Value = "Zapallo"
$ws.Range("H403").Value = "Camote"
$ws.Range("I403").Value = "1a (cosecha)"
$ws.Range("J403").Value = 210
$ws.Range("K403").Value = 500
$ws.Range("L403").Value = 550
$ws.Range("M403").Value = 524
$ws.Range("N403").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O403").Value = "Provincia de Talca"
$ws.Range("P403").Value = 524
$ws.Range("Q403").Value = 1
$ws.Range("R403").Value = "Hortaliza"

# --- New row 404: Paine, 1a (cosecha) ---
$ws.Range("A404").Value = 3
$ws.Range("B404").Value = "Femacal de La Calera"
$ws.Range("C404").Value = "Coquimbo"
$ws.Range("D404").Value = 44585
$ws.Range("E404").Value = 5
$ws.Range("F404").Value = 100112045
$ws.Range("G404").Value = "Zapallo"
$ws.Range("H404").Value = "Paine"
$ws.Range("I404").Value = "1a (cosecha)"
$ws.Range("J404").Value = 170
$ws.Range("K404").Value = 250
$ws.Range("L404").Value = 300
$ws.Range("M404").Value = 274
$ws.Range("N404").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O404").Value = "Provincia de Talca"
$ws.Range("P404").Value = 274
$ws.Range("Q404").Value = 1
$ws.Range("R404").Value = "Hortaliza"
